# Apply "repull data, push all data, mean calculation" changes:
# Updates the dSF column (F) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -6
$ws.Range("F4").Value  = -10
$ws.Range("F6").Value  = -5
$ws.Range("F7").Value  = -5
$ws.Range("F9").Value  = -16
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = 10
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -9
